$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the existing data range into a native Excel Table ("ListObject"),
# matching the author's "Format as Table" action (default TableStyleMedium2 style).
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:E101"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "VehicleInventory"
$tbl.TableStyle = "TableStyleMedium2"

# Widen column B slightly (author resized it while reviewing the new table).
$ws.Columns.Item(2).ColumnWidth = 8.938

# Select the whole table body (mirrors the author's final selection state)
# and make sure the view is scrolled back to the top of the sheet.
$ws.Range("A1").Select() | Out-Null
$ws.Range("A1:E101").Select() | Out-Null
